# "plot with real data for inference"
#
# Updates the inference-time placeholder numbers (NWS external-page-size
# row and SS row) on the "timeoverhead" sheet with the real measured
# values, and restores the selections left behind on the
# "timeoverhead_backup" and "timeoverhead" sheets after the edit.

$wb = $excel.ActiveWorkbook

$wsBackup = $wb.Worksheets.Item("timeoverhead_backup")
$wsMain   = $wb.Worksheets.Item("timeoverhead")

# --- Real data for the "NWS - External page size" row (row 7) ---
$wsMain.Range("B7").Value = 30.26
$wsMain.Range("C7").Value = 40.58
$wsMain.Range("D7").Value = 40.58
$wsMain.Range("E7").Value = 41.06
$wsMain.Range("F7").Value = 19.72

# --- Real data for the "SS" row (row 8) ---
$wsMain.Range("B8").Value = 18.21
$wsMain.Range("C8").Value = 17.73
$wsMain.Range("D8").Value = 17.73
$wsMain.Range("E8").Value = 18.09
$wsMain.Range("F8").Value = 8.81

# --- Leave the selection on timeoverhead_backup where the author left it ---
$wsBackup.Activate()
$wsBackup.Range("A1").Select()
$wsBackup.Range("H27").Select()

# --- Re-activate timeoverhead (it stays the tab in view) and park the
#     selection on the cell the author ended up on ---
$wsMain.Activate()
$wsMain.Range("E17").Select()
